$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.294.88"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.62"
$ws.Range("E3").Value = "  -6.34%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.12"
$ws.Range("E5").Value = "  -4.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.01"
$ws.Range("E6").Value = "  -4.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.373.68"
$ws.Range("E9").Value = "  -5.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0953"
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.77"
$ws.Range("E12").Value = "  -8.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.315"
$ws.Range("E13").Value = "  -5.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.798.10"
$ws.Range("E14").Value = "  -4.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.317.23"
$ws.Range("E15").Value = "  -3.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.40"
$ws.Range("E16").Value = "  -4.15%  "
$ws.Range("E17").Value = "  -4.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.373.74"
$ws.Range("E18").Value = "  -5.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.26"
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.05"
$ws.Range("E20").Value = "  -4.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.36"
$ws.Range("E21").Value = "  -3.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.86"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.389"
$ws.Range("E26").Value = "  -5.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.467.01"
$ws.Range("E27").Value = "  -5.96%  "
$ws.Range("E28").Value = "  -5.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  -3.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.27"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.67"
$ws.Range("E31").Value = "  -3.97%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0711"
$ws.Range("E33").Value = "  -6.64%  "
$ws.Range("E34").Value = "  -7.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.60"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("E38").Value = "  -4.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.68"
$ws.Range("E39").Value = "  -6.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.58"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.787"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  -6.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.74"
$ws.Range("E44").Value = "  -3.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.77"
$ws.Range("E45").Value = "  -7.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "254.83"
$ws.Range("E46").Value = "  -7.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.567"
$ws.Range("E47").Value = "  -4.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0899"
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0484"
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("E50").Value = "  -6.20%  "
$ws.Range("E51").Value = "  -6.35%  "
